{"js": "// Replace each arithmetic-expression cell text in the table with its new value.\n// Each [oldText, newText] pair corresponds to one <w:t> run changed by the diff,\n// in document order. All oldText values are unique within the document, so an\n// exact-text search safely targets the correct run.\nconst pairs = [\n  [\"12-8=\", \"78-67=\"],\n  [\"23+53=\", \"22+32=\"],\n  [\"84-17=\", \"13+54=\"],\n  [\"7+82=\", \"81-11=\"],\n  [\"95-0=\", \"79+3=\"],\n  [\"95-31=\", \"96-31=\"],\n  [\"28+65=\", \"46+20=\"],\n  [\"93+0=\", \"93-56=\"],\n  [\"79-27=\", \"51-14=\"],\n  [\"70-23=\", \"26+0=\"],\n  [\"98-0=\", \"2+51=\"],\n  [\"12+77=\", \"69-66=\"],\n  [\"48+26=\", \"25+1=\"],\n  [\"27+1=\", \"85-54=\"],\n  [\"15+13=\", \"29+32=\"],\n  [\"20+31=\", \"5+66=\"],\n  [\"25+74=\", \"29+16=\"],\n  [\"49+38=\", \"58-54=\"],\n  [\"31-25=\", \"31+5=\"],\n  [\"96-69=\", \"70+14=\"],\n  [\"83+11=\", \"8+28=\"],\n  [\"8+67=\", \"15+2=\"],\n  [\"37-32=\", \"30-27=\"],\n  [\"55-31=\", \"49+0=\"],\n  [\"85-70=\", \"93+6=\"],\n  [\"18-11=\", \"16-6=\"],\n  [\"0+64=\", \"91-69=\"],\n  [\"39+37=\", \"92-80=\"],\n  [\"81-49=\", \"12+13=\"],\n  [\"53-26=\", \"48+30=\"],\n  [\"27-9=\", \"39-25=\"],\n  [\"52+35=\", \"18+55=\"],\n  [\"22+20=\", \"84-42=\"],\n  [\"48+21=\", \"70-10=\"],\n  [\"65-12=\", \"55-41=\"],\n  [\"8+8=\", \"5+83=\"],\n  [\"13+63=\", \"78+4=\"],\n  [\"6-2=\", \"15+36=\"],\n  [\"77+12=\", \"89-24=\"],\n  [\"23+67=\", \"88-7=\"],\n  [\"54-5=\", \"6+82=\"],\n  [\"82-58=\", \"13+46=\"],\n  [\"53-29=\", \"89-69=\"],\n  [\"52-1=\", \"38+50=\"],\n  [\"38+46=\", \"76-26=\"],\n  [\"76-37=\", \"38+54=\"],\n  [\"38-28=\", \"92+1=\"],\n  [\"81+2=\", \"12-11=\"],\n  [\"51+34=\", \"38-6=\"],\n  [\"35+48=\", \"42-35=\"],\n  [\"85-31=\", \"33-25=\"],\n  [\"70+4=\", \"81-22=\"],\n  [\"76-64=\", \"13+8=\"],\n  [\"91-25=\", \"17+20=\"],\n  [\"20+27=\", \"90-35=\"],\n  [\"73-56=\", \"29-1=\"],\n  [\"56-28=\", \"60-53=\"],\n  [\"91-32=\", \"30-20=\"],\n  [\"28-4=\", \"54-39=\"],\n  [\"43-28=\", \"76-76=\"],\n  [\"63+8=\", \"52+16=\"],\n  [\"66-43=\", \"74-35=\"],\n  [\"11+54=\", \"33+24=\"],\n  [\"38-36=\", \"20+26=\"],\n  [\"9+2=\", \"72-28=\"],\n  [\"54-1=\", \"57+2=\"],\n  [\"52-40=\", \"77-52=\"],\n  [\"52-14=\", \"2+13=\"],\n  [\"64+28=\", \"43+42=\"],\n  [\"73-10=\", \"59-34=\"],\n  [\"66-38=\", \"59+27=\"],\n  [\"11+72=\", \"7+76=\"],\n  [\"49-7=\", \"92-43=\"],\n  [\"89-25=\", \"38+33=\"],\n  [\"9-4=\", \"92-13=\"],\n  [\"56+2=\", \"98-78=\"],\n  [\"97-2=\", \"34-9=\"],\n  [\"51+0=\", \"91-63=\"],\n  [\"3+66=\", \"77+18=\"],\n  [\"73-69=\", \"89+5=\"],\n  [\"15+60=\", \"23+6=\"],\n  [\"92-72=\", \"44+41=\"],\n  [\"1+18=\", \"75-65=\"],\n  [\"39-38=\", \"95-68=\"],\n  [\"24+8=\", \"64+12=\"],\n  [\"74-29=\", \"43+17=\"],\n  [\"2-0=\", \"56+3=\"],\n  [\"68-30=\", \"62-14=\"],\n  [\"26+70=\", \"65-11=\"],\n  [\"5+39=\", \"24+44=\"],\n  [\"31+45=\", \"17+74=\"],\n  [\"57-36=\", \"67+6=\"],\n  [\"71-44=\", \"62-1=\"],\n  [\"42-38=\", \"65-24=\"],\n  [\"94-71=\", \"86-64=\"],\n  [\"91-61=\", \"34-22=\"],\n  [\"34+37=\", \"9+61=\"],\n  [\"60+28=\", \"97-22=\"],\n  [\"48+17=\", \"33+14=\"],\n  [\"81-5=\", \"5+8=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  // Replace the matching run in place (texts are unique, so items[0] is\n  // unambiguous).\n  results.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Replace each arithmetic-expression cell's text in the table with its new\n# value. Each pair is [oldText, newText] in document order, matching the\n# <w:t> run changes in the diff. All oldText values are unique within the\n# document, so Find/Replace with whole-word + case-sensitive matching\n# unambiguously targets the correct cell; wdReplaceAll is used defensively\n# but only ever matches a single occurrence per pair.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"12-8=\", \"78-67=\"),\n  @(\"23+53=\", \"22+32=\"),\n  @(\"84-17=\", \"13+54=\"),\n  @(\"7+82=\", \"81-11=\"),\n  @(\"95-0=\", \"79+3=\"),\n  @(\"95-31=\", \"96-31=\"),\n  @(\"28+65=\", \"46+20=\"),\n  @(\"93+0=\", \"93-56=\"),\n  @(\"79-27=\", \"51-14=\"),\n  @(\"70-23=\", \"26+0=\"),\n  @(\"98-0=\", \"2+51=\"),\n  @(\"12+77=\", \"69-66=\"),\n  @(\"48+26=\", \"25+1=\"),\n  @(\"27+1=\", \"85-54=\"),\n  @(\"15+13=\", \"29+32=\"),\n  @(\"20+31=\", \"5+66=\"),\n  @(\"25+74=\", \"29+16=\"),\n  @(\"49+38=\", \"58-54=\"),\n  @(\"31-25=\", \"31+5=\"),\n  @(\"96-69=\", \"70+14=\"),\n  @(\"83+11=\", \"8+28=\"),\n  @(\"8+67=\", \"15+2=\"),\n  @(\"37-32=\", \"30-27=\"),\n  @(\"55-31=\", \"49+0=\"),\n  @(\"85-70=\", \"93+6=\"),\n  @(\"18-11=\", \"16-6=\"),\n  @(\"0+64=\", \"91-69=\"),\n  @(\"39+37=\", \"92-80=\"),\n  @(\"81-49=\", \"12+13=\"),\n  @(\"53-26=\", \"48+30=\"),\n  @(\"27-9=\", \"39-25=\"),\n  @(\"52+35=\", \"18+55=\"),\n  @(\"22+20=\", \"84-42=\"),\n  @(\"48+21=\", \"70-10=\"),\n  @(\"65-12=\", \"55-41=\"),\n  @(\"8+8=\", \"5+83=\"),\n  @(\"13+63=\", \"78+4=\"),\n  @(\"6-2=\", \"15+36=\"),\n  @(\"77+12=\", \"89-24=\"),\n  @(\"23+67=\", \"88-7=\"),\n  @(\"54-5=\", \"6+82=\"),\n  @(\"82-58=\", \"13+46=\"),\n  @(\"53-29=\", \"89-69=\"),\n  @(\"52-1=\", \"38+50=\"),\n  @(\"38+46=\", \"76-26=\"),\n  @(\"76-37=\", \"38+54=\"),\n  @(\"38-28=\", \"92+1=\"),\n  @(\"81+2=\", \"12-11=\"),\n  @(\"51+34=\", \"38-6=\"),\n  @(\"35+48=\", \"42-35=\"),\n  @(\"85-31=\", \"33-25=\"),\n  @(\"70+4=\", \"81-22=\"),\n  @(\"76-64=\", \"13+8=\"),\n  @(\"91-25=\", \"17+20=\"),\n  @(\"20+27=\", \"90-35=\"),\n  @(\"73-56=\", \"29-1=\"),\n  @(\"56-28=\", \"60-53=\"),\n  @(\"91-32=\", \"30-20=\"),\n  @(\"28-4=\", \"54-39=\"),\n  @(\"43-28=\", \"76-76=\"),\n  @(\"63+8=\", \"52+16=\"),\n  @(\"66-43=\", \"74-35=\"),\n  @(\"11+54=\", \"33+24=\"),\n  @(\"38-36=\", \"20+26=\"),\n  @(\"9+2=\", \"72-28=\"),\n  @(\"54-1=\", \"57+2=\"),\n  @(\"52-40=\", \"77-52=\"),\n  @(\"52-14=\", \"2+13=\"),\n  @(\"64+28=\", \"43+42=\"),\n  @(\"73-10=\", \"59-34=\"),\n  @(\"66-38=\", \"59+27=\"),\n  @(\"11+72=\", \"7+76=\"),\n  @(\"49-7=\", \"92-43=\"),\n  @(\"89-25=\", \"38+33=\"),\n  @(\"9-4=\", \"92-13=\"),\n  @(\"56+2=\", \"98-78=\"),\n  @(\"97-2=\", \"34-9=\"),\n  @(\"51+0=\", \"91-63=\"),\n  @(\"3+66=\", \"77+18=\"),\n  @(\"73-69=\", \"89+5=\"),\n  @(\"15+60=\", \"23+6=\"),\n  @(\"92-72=\", \"44+41=\"),\n  @(\"1+18=\", \"75-65=\"),\n  @(\"39-38=\", \"95-68=\"),\n  @(\"24+8=\", \"64+12=\"),\n  @(\"74-29=\", \"43+17=\"),\n  @(\"2-0=\", \"56+3=\"),\n  @(\"68-30=\", \"62-14=\"),\n  @(\"26+70=\", \"65-11=\"),\n  @(\"5+39=\", \"24+44=\"),\n  @(\"31+45=\", \"17+74=\"),\n  @(\"57-36=\", \"67+6=\"),\n  @(\"71-44=\", \"62-1=\"),\n  @(\"42-38=\", \"65-24=\"),\n  @(\"94-71=\", \"86-64=\"),\n  @(\"91-61=\", \"34-22=\"),\n  @(\"34+37=\", \"9+61=\"),\n  @(\"60+28=\", \"97-22=\"),\n  @(\"48+17=\", \"33+14=\"),\n  @(\"81-5=\", \"5+8=\")\n)\n\nforeach ($pair in $pairs) {\n  $findText = $pair[0]\n  $replaceText = $pair[1]\n  $range = $d.Content\n  $range.Find.ClearFormatting()\n  $range.Find.Replacement.ClearFormatting()\n  $found = $range.Find.Execute(\n      $findText,    # FindText\n      $true,        # MatchCase\n      $true,        # MatchWholeWord\n      $false,       # MatchWildcards\n      $false,       # MatchSoundsLike\n      $false,       # MatchAllWordForms\n      $true,        # Forward\n      1,            # Wrap (wdFindContinue)\n      $false,       # Format\n      $replaceText, # ReplaceWith\n      2             # Replace (wdReplaceAll)\n  )\n  if (-not $found) {\n    throw \"No match found for: $findText\"\n  }\n}\n\n"}
